# Update "想去人数" (want-to-go count) figures across the four sheets to
# reflect the newly scraped numbers recorded in commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 138
$ws.Range("F4").Value = 40
$ws.Range("F6").Value = 279
$ws.Range("F7").Value = 13211
$ws.Range("F8").Value = 72
$ws.Range("F10").Value = 297
$ws.Range("F11").Value = 4739
$ws.Range("F13").Value = 3557
$ws.Range("F17").Value = 180
$ws.Range("F21").Value = 77
$ws.Range("F25").Value = 4246
$ws.Range("F26").Value = 428
$ws.Range("F27").Value = 1937
$ws.Range("F29").Value = 259
$ws.Range("F30").Value = 6990
$ws.Range("F31").Value = 20
$ws.Range("F34").Value = 2056
$ws.Range("F36").Value = 116
$ws.Range("F40").Value = 232
$ws.Range("F43").Value = 7
$ws.Range("F45").Value = 1240
$ws.Range("F46").Value = 1854
$ws.Range("F47").Value = 78
$ws.Range("F48").Value = 172
$ws.Range("F49").Value = 1188

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 131

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 480
$ws.Range("F3").Value = 650
$ws.Range("F4").Value = 36

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 138
$ws.Range("F4").Value = 40
$ws.Range("F5").Value = 480
$ws.Range("F6").Value = 650
$ws.Range("F7").Value = 279
$ws.Range("F8").Value = 13212
$ws.Range("F10").Value = 297
$ws.Range("F11").Value = 4739
$ws.Range("F12").Value = 3557
$ws.Range("F15").Value = 180
$ws.Range("F18").Value = 77
$ws.Range("F23").Value = 4246
$ws.Range("F24").Value = 428
$ws.Range("F25").Value = 1937
$ws.Range("F27").Value = 259
$ws.Range("F28").Value = 6990
$ws.Range("F30").Value = 20
$ws.Range("F33").Value = 2056
$ws.Range("F35").Value = 116
$ws.Range("F38").Value = 232
$ws.Range("F41").Value = 7
$ws.Range("F44").Value = 1240
$ws.Range("F45").Value = 1854
$ws.Range("F46").Value = 78
$ws.Range("F48").Value = 172
$ws.Range("F49").Value = 1188
